$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "28h 50m"
$ws.Range("B4").Value = "56h 30m"

$ws.Range("B4").Select()
